$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 15 ("Manual - Team FHIR When Ready.pdf"),
# shifting it and everything below it down by one row.
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the "Special Instructions" catalog entry,
# matching the same layout used by the other top-level "2 / file name" rows
# (e.g. old row 14 / new row 14: B=2, C=<file name>).
$newRow = $ws.Rows.Item(15)
$bCell = $ws.Cells.Item(15, 2)

# Writing a number directly into the freshly-inserted (quote-prefixed) cell
# would silently reformat it to a non-quote-prefixed style, so reset the
# cell's format first, assign the value, then copy the real formatting back
# from the row above (which already carries the desired style).
$bCell.ClearFormats()
$bCell.Value = 2
$ws.Cells.Item(15, 3).Value = "Special Instructions - Team FHIR When Ready.pdf"

$ws.Cells.Item(14, 2).Copy()
$bCell.PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match the taller row height used by the other "2 / file name" rows
# (e.g. row 14), which use a 14pt Wingdings bullet font.
$ws.Rows.Item(15).RowHeight = 18

# Reflect the new last used row in the print area.
$ws.PageSetup.PrintArea = '$A$1:$P$44'

# Match the saved selection state from the edit.
$ws.Range("D24").Select()

Write-Output "applied catalog edit"
